$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "LP1912": refresh timestamp/header + new set of 8 data rows
# (previously 5 rows, now 8 -> rows 6..13)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 03:39:35"
$ws1.Range("A3").Value = "Total filas: 8"

$sheet1Data = @(
    @("03:39:35", "03:41", "14_ABASTO", 2, "LP1912"),
    @("03:39:35", "04:01", "81_EL PELIGRO", 22, "LP1912"),
    @("03:39:35", "04:46", "215A_EL PATO", 67, "LP1912"),
    @("03:39:35", "04:53", "11_ETCHEVERRY", 74, "LP1912"),
    @("03:39:35", "05:16", "17_ROMERO", 97, "LP1912"),
    @("03:39:35", "05:22", "23_HERNANDEZ", 103, "LP1912"),
    @("03:39:35", "05:31", "14_ABASTO", 112, "LP1912"),
    @("03:39:35", "05:34", "215B_EL PATO", 115, "LP1912")
)

$row = 6
foreach ($rec in $sheet1Data) {
    $ws1.Cells.Item($row, 1).Value = $rec[0]
    $ws1.Cells.Item($row, 2).Value = $rec[1]
    $ws1.Cells.Item($row, 3).Value = $rec[2]
    $ws1.Cells.Item($row, 4).Value = $rec[3]
    $ws1.Cells.Item($row, 5).Value = $rec[4]
    $row++
}

# ---------------------------------------------------------------
# Sheet "LP1912-215": refresh timestamp + update its 2 data rows
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 03:39:35"

$sheet2Data = @(
    @("03:39:35", "04:46", "215A_EL PATO", 67, "LP1912"),
    @("03:39:35", "05:34", "215B_EL PATO", 115, "LP1912")
)

$row = 6
foreach ($rec in $sheet2Data) {
    $ws2.Cells.Item($row, 1).Value = $rec[0]
    $ws2.Cells.Item($row, 2).Value = $rec[1]
    $ws2.Cells.Item($row, 3).Value = $rec[2]
    $ws2.Cells.Item($row, 4).Value = $rec[3]
    $ws2.Cells.Item($row, 5).Value = $rec[4]
    $row++
}

# ---------------------------------------------------------------
# Sheet "6203-6173": refresh timestamp only (no data rows)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 03:39:35"
